$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q3" right after "总计" (i.e. before
#    the sheet currently in position 2, which is "2022-Q2"). This shifts
#    "2022-Q2" .. "2021-Q1" one slot to the right, matching the diff.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# Match the sheet-level print/outline defaults used by every other sheet in
# this workbook (summaryBelow="1" summaryRight="1", 0.75"/1"/0.5" margins).
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 0.75 * 72
$newSheet.PageSetup.RightMargin = 0.75 * 72
$newSheet.PageSetup.TopMargin = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the fund holdings table.
#    Columns B-G carry text-formatted values (fund codes / percentages keep
#    their original string formatting, incl. leading zeros / trailing
#    zeros), while columns A and H are plain numbers.
# ---------------------------------------------------------------------------
$newSheet.Range("B1:H1").NumberFormat = "@"
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("B2:G4").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "167506"
$newSheet.Range("C2").Value = "安信深圳科技指数（LOF）A"
$newSheet.Range("D2").Value = "0.82"
$newSheet.Range("E2").Value = "93.20"
$newSheet.Range("F2").Value = "5.25"
$newSheet.Range("G2").Value = "0.0430"
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "167507"
$newSheet.Range("C3").Value = "安信深圳科技指数（LOF）C"
$newSheet.Range("D3").Value = "0.30"
$newSheet.Range("E3").Value = "93.20"
$newSheet.Range("F3").Value = "5.25"
$newSheet.Range("G3").Value = "0.0158"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "159932"
$newSheet.Range("C4").Value = "大成中证500深市ETF"
$newSheet.Range("D4").Value = "0.32"
$newSheet.Range("E4").Value = "95.99"
$newSheet.Range("F4").Value = "1.04"
$newSheet.Range("G4").Value = "0.0033"
$newSheet.Range("H4").Value = 5

# Match the header / index-column style used by every other sheet in this
# workbook (bold, centered, thin border = style index 2 on the template
# sheets) by copying formats from the analogous cells on the "总计" sheet,
# then drop the temporary text-number-format from the plain data cells so
# they end up with no explicit style, just like the other quarter sheets.
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$newSheet.Range("B2:G4").ClearFormats()

# ---------------------------------------------------------------------------
# 3. Insert the matching summary row at the top of the "总计" sheet's data
#    (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.06

# Restore the original active sheet/selection ("总计", cell A1) so the
# workbook-level view state is unaffected by our edits.
[void]$totalSheet.Activate()
[void]$totalSheet.Range("A1").Select()
